$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.726.62"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.45%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.380.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.83%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.12%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'561.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.16%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'175.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.24%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.09%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.374.94"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.79%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.00%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +1.10%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -0.40%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'53.64"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.49%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  -1.16%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +0.12%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.924.66"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.50%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "'Chainlink"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'18.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.04%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "'TRON"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'0.119"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.09%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "'WrappedEther"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'3.374.32"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.79%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'65.653.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.34%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'11.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.81%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  -0.21%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'465.87"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.69%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'4.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -2.38%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'89.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +2.99%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'14.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +6.06%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -1.72%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +0.00%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'10.58"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -3.28%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'8.69"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -2.29%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'31.03"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.79%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -2.68%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'11.41"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.29%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'579.42"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.71%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'61.92"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.05%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -0.84%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -0.04%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +2.44%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +0.97%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'35.91"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.08%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.374"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.02%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.0₃0739"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -2.95%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'3.094.33"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.07%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -1.19%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -0.74%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -1.38%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'Fetch.AI"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'2.44"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -2.31%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "'ApeXProtocol"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'3.16"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.18%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +0.18%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'140.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +1.95%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -1.80%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'8.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.52%  "
$ws.Range("E51").Style = "Normal"
